$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# LIKING_PROMPT / DIFFICULTY_PROMPT rows (18-19)
$ws.Range("B18").Value = "Das Stück hat mir gefallen."
$ws.Range("C18").Value = "I liked the piece."
$ws.Range("B19").Value = "Die Aufgabe war schwierig."

# CREDITS row (37) - swap &copy; spacing, English then German
$ws.Range("C37").Value = "John Knowles Paine: Symphony No. 1 in C minor, Op. 23: Allegro con brio. From the sound recording New World Records #80374-2.  <br/> &copy;1989 Anthology of Recorded Music, Inc. Used by permission."
$ws.Range("B37").Value = "John Knowles Paine: Symphonie Nr. 1 in c-Moll, Op. 23: Allegro con brio. Aus der Aufnahme New World Records #80374-2. <br/> &copy;1989 Anthology of Recorded Music, Inc. Gebrauch mit freundlicher Genehmigung."

# CONTINUE_MAIN_TEST row (22) - new wording with <br> and "have fun" addendum
$ws.Range("B22").Value = "Beginne mit dem Experiment. <br> Viel Vergnügen!"
$ws.Range("C22").Value = "Begin the experiment. <br> Have fun!"

# NUM_LIKERT1 / NUM_LIKERT6 rows (29 / 34) - reworded numeric anchor labels
$ws.Range("B29").Value = "Trifft gar nicht zu 1"
$ws.Range("B34").Value = "6 Trifft sehr zu"
$ws.Range("C34").Value = "6 Completely agree"
$ws.Range("C29").Value = "Completely disagree 1"

# Update selection / scroll position to match the new view (was topLeftCell A22 / C37)
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
$ws.Range("C29").Select()
